# Update NATMI ligand-receptor pair metrics (Cdh1-Egfr) with values
# recomputed from the new TPM input data (commit: "update scripts wuth new tpm").
# Columns E-J depend on the sending cluster, M-P on the target cluster, and
# Q-T on the sending/target edge; A-D, K and L are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1443736666666667
$ws.Range("H2").Value = 0.433121
$ws.Range("I2").Value = 0.7378778224885942
$ws.Range("J2").Value = 0.7378778224885942
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 0.1979184427898889
$ws.Range("R2").Value = 1.781265985109
$ws.Range("S2").Value = 0.008139259527623227
$ws.Range("T2").Value = 0.008139259527623226
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1443736666666667
$ws.Range("H3").Value = 0.433121
$ws.Range("I3").Value = 0.7378778224885942
$ws.Range("J3").Value = 0.7378778224885942
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 13.41429506565922
$ws.Range("R3").Value = 120.728655590933
$ws.Range("S3").Value = 0.5516536376320657
$ws.Range("T3").Value = 0.5516536376320655
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1443736666666667
$ws.Range("H4").Value = 0.433121
$ws.Range("I4").Value = 0.7378778224885942
$ws.Range("J4").Value = 0.7378778224885942
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 4.290566502854667
$ws.Range("R4").Value = 38.615098525692
$ws.Range("S4").Value = 0.1764465897922121
$ws.Range("T4").Value = 0.176446589792212
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1443736666666667
$ws.Range("H5").Value = 0.433121
$ws.Range("I5").Value = 0.7378778224885942
$ws.Range("J5").Value = 0.7378778224885942
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.03983861395366667
$ws.Range("R5").Value = 0.358547525583
$ws.Range("S5").Value = 0.001638335536693352
$ws.Range("T5").Value = 0.001638335536693351
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.051287
$ws.Range("H6").Value = 0.153861
$ws.Range("I6").Value = 0.2621221775114058
$ws.Range("J6").Value = 0.2621221775114058
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 0.07030813450766667
$ws.Range("R6").Value = 0.632773210569
$ws.Range("S6").Value = 0.002891373565769467
$ws.Range("T6").Value = 0.002891373565769467
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.051287
$ws.Range("H7").Value = 0.153861
$ws.Range("I7").Value = 0.2621221775114058
$ws.Range("J7").Value = 0.2621221775114058
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("Q7").Value = 4.765266179883667
$ws.Range("R7").Value = 42.887395618953
$ws.Range("S7").Value = 0.1959682867829249
$ws.Range("T7").Value = 0.1959682867829249
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.051287
$ws.Range("H8").Value = 0.153861
$ws.Range("I8").Value = 0.2621221775114058
$ws.Range("J8").Value = 0.2621221775114058
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 1.524171888908
$ws.Range("R8").Value = 13.717547000172
$ws.Range("S8").Value = 0.06268051826630326
$ws.Range("T8").Value = 0.06268051826630326
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.051287
$ws.Range("H9").Value = 0.153861
$ws.Range("I9").Value = 0.2621221775114058
$ws.Range("J9").Value = 0.2621221775114058
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 0.014152186067
$ws.Range("R9").Value = 0.127369674603
$ws.Range("S9").Value = 0.000581998896408107
$ws.Range("T9").Value = 0.000581998896408107
